$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# "About" sheet
# ---------------------------------------------------------------------------
$about = $wb.Worksheets.Item("About")

# ---------------------------------------------------------------------------
# "Data" sheet
# ---------------------------------------------------------------------------
$data = $wb.Worksheets.Item("Data")

$data.Range("A1").Value = "Excerpt from Table 6-2:"

# Updated "last edited" date
$about.Range("C1").Value = 44515

# Source block
$about.Range("B3").Value = "US EPA"
$about.Range("B4").Value = 2021
$about.Range("B5").Value = "Draft Inventory of US Greenhouse Gas Emissions Emissions and Sinks"
$about.Range("B6").Value = "https://www.epa.gov/sites/production/files/2021-02/documents/us-ghg-inventory-2021-main-text.pdf"
$about.Range("B7").Value = "Table 6-3"

# New "California" label next to the title
$about.Range("B1").Value = "California"

# Years
$data.Range("B3").Value = 2015
$data.Range("C3").Value = 2016
$data.Range("D3").Value = 2017
$data.Range("E3").Value = 2018
$data.Range("F3").Value = 2019

# CO2
$data.Range("B4").Value = -791695
$data.Range("C4").Value = -855998
$data.Range("D4").Value = -792046
$data.Range("E4").Value = -824885
$data.Range("F4").Value = -812695

# CH4
$data.Range("B5").Value = 663
$data.Range("C5").Value = 308
$data.Range("D5").Value = 614
$data.Range("E5").Value = 552
$data.Range("F5").Value = 552

# N2O
$data.Range("B6").Value = 38
$data.Range("C6").Value = 18
$data.Range("D6").Value = 36
$data.Range("E6").Value = 32
$data.Range("F6").Value = 32

# ---------------------------------------------------------------------------
# Sheet view / selection updates
# ---------------------------------------------------------------------------
$about.Range("B8").Select()

$data.Range("F6").Select()

$rpe = $wb.Worksheets.Item("RPEpUACE")
$rpe.Range("B2:B13").Select()

# Re-activate "About" so it stays the selected tab, matching the source file
$about.Activate()
